$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "semantic_aspect_model_schema": header row field-name fixups
# (single underscore separator -> double underscore separator) and
# two column-width tweaks (col H / 8 and col J / 10).
# -----------------------------------------------------------------
$wsSchema = $wb.Worksheets.Item("semantic_aspect_model_schema")

$wsSchema.Range("D1").Value = "criticalRawMaterials[0]__criticalRawMaterialName"
$wsSchema.Range("E1").Value = "criticalRawMaterials[0]__percentageOfMaterialWeight"
$wsSchema.Range("F1").Value = "hazardousMaterials[0]__hazardClassification"
$wsSchema.Range("G1").Value = "hazardousMaterials[0]__locationOfHazardousSubstances"
$wsSchema.Range("H1").Value = "hazardousMaterials[0]__enviromentalImpact"
$wsSchema.Range("I1").Value = "hazardousMaterials[0]__hazardousMaterialCas"
$wsSchema.Range("J1").Value = "hazardousMaterials[0]__clpIndex"
$wsSchema.Range("K1").Value = "hazardousMaterials[0]__rangeScipConcentration"

# width="48" -> width="49.2" (column H, index 8)
$wsSchema.Columns.Item(8).ColumnWidth = 48.3
# width="36" -> width="37.2" (column J, index 10)
$wsSchema.Columns.Item(10).ColumnWidth = 36.3

# -----------------------------------------------------------------
# Sheet "description": legend wording + field-name fixups
# -----------------------------------------------------------------
$wsDescription = $wb.Worksheets.Item("description")

$wsDescription.Range("A3").Value = "1. Columns highlighted in olive green are digital twin fields."

$wsDescription.Range("B5").Value = "Digital Twin Field Name: id"
$wsDescription.Range("B6").Value = "Digital Twin Field Name: manufacturerPartId"
$wsDescription.Range("B7").Value = "Digital Twin Field Name: partInstanceId"

$wsDescription.Range("A8").Value = "criticalRawMaterials[0]__criticalRawMaterialName"
$wsDescription.Range("A9").Value = "criticalRawMaterials[0]__percentageOfMaterialWeight"
$wsDescription.Range("A10").Value = "hazardousMaterials[0]__hazardClassification"
$wsDescription.Range("A11").Value = "hazardousMaterials[0]__locationOfHazardousSubstances"
$wsDescription.Range("A12").Value = "hazardousMaterials[0]__enviromentalImpact"
$wsDescription.Range("A13").Value = "hazardousMaterials[0]__hazardousMaterialCas"
$wsDescription.Range("A14").Value = "hazardousMaterials[0]__clpIndex"
$wsDescription.Range("A15").Value = "hazardousMaterials[0]__rangeScipConcentration"
